$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B:F columns (ME, MAE, MSE, RMSE, SE) for rows 2-11
$ws.Range("B2").Value = 0.09104566687622702
$ws.Range("C2").Value = 0.8305735710435281
$ws.Range("D2").Value = 1.571787696088291
$ws.Range("E2").Value = 1.253709574059435
$ws.Range("F2").Value = 1.262598595620691

$ws.Range("B3").Value = 0.2881368179803065
$ws.Range("C3").Value = 0.9230054831243015
$ws.Range("D3").Value = 1.736110936727699
$ws.Range("E3").Value = 1.317615625562971
$ws.Range("F3").Value = 1.298518337594399

$ws.Range("B4").Value = 0.2572367959422796
$ws.Range("C4").Value = 0.9122950025140623
$ws.Range("D4").Value = 1.616069083570048
$ws.Range("E4").Value = 1.271247058431227
$ws.Range("F4").Value = 1.257588522620955

$ws.Range("B5").Value = 0.2382191396178211
$ws.Range("C5").Value = 0.7774418406714947
$ws.Range("D5").Value = 1.214992806219648
$ws.Range("E5").Value = 1.102267121082566
$ws.Range("F5").Value = 1.089587334238462
$ws.Range("G5").Value = 41

$ws.Range("B6").Value = 0.3461003139263453
$ws.Range("C6").Value = 0.7692343642516598
$ws.Range("D6").Value = 1.083647203965686
$ws.Range("E6").Value = 1.040983767388179
$ws.Range("F6").Value = 0.9979932377299109
$ws.Range("G6").Value = 31

$ws.Range("B7").Value = 0.390978884314848
$ws.Range("C7").Value = 0.7686969224094469
$ws.Range("D7").Value = 1.1171720694664
$ws.Range("E7").Value = 1.056963608392645
$ws.Range("F7").Value = 0.9993733726530692
$ws.Range("G7").Value = 29

$ws.Range("B8").Value = 0.3395756546653472
$ws.Range("C8").Value = 0.7550656853727538
$ws.Range("D8").Value = 1.132181430349046
$ws.Range("E8").Value = 1.064040145083373
$ws.Range("F8").Value = 1.027608963675956
$ws.Range("G8").Value = 27

$ws.Range("B9").Value = 0.3539796401532963
$ws.Range("C9").Value = 0.862300489931706
$ws.Range("D9").Value = 1.388483373928287
$ws.Range("E9").Value = 1.178339243990578
$ws.Range("F9").Value = 1.154711459326994
$ws.Range("G9").Value = 19

$ws.Range("B10").Value = -0.1625119703050184
$ws.Range("C10").Value = 0.4241755480675304
$ws.Range("D10").Value = 0.2902502803524973
$ws.Range("E10").Value = 0.5387488100706092
$ws.Range("F10").Value = 0.5364938090230728
$ws.Range("G10").Value = 12

$ws.Range("B11").Value = -0.2103864081339538
$ws.Range("C11").Value = 0.4484001654209663
$ws.Range("D11").Value = 0.297018969647196
$ws.Range("E11").Value = 0.5449944675381539
$ws.Range("F11").Value = 0.5620904385858309
